$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: convert A3 and Q3:U3 from text to numeric values
$ws.Range("A3").Value = 71
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 2

# New rows 4-8: fully numeric/text typed data
# Row 4
$ws.Range("A4").Value = 36
$ws.Range("B4").Value = "Masculino"
$ws.Range("C4").Value = "sim"
$ws.Range("D4").Value = "sim"
$ws.Range("E4").Value = "nao"
$ws.Range("F4").Value = "sim"
$ws.Range("G4").Value = "nao"
$ws.Range("H4").Value = "sim"
$ws.Range("I4").Value = "nao"
$ws.Range("J4").Value = "nao"
$ws.Range("K4").Value = "nao"
$ws.Range("L4").Value = "sim"
$ws.Range("M4").Value = "sim"
$ws.Range("N4").Value = "nao"
$ws.Range("O4").Value = "sim"
$ws.Range("P4").Value = "sim"
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 3
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = "bom"
$ws.Range("W4").Value = "reprovado"
$ws.Range("X4").Value = "'"
$ws.Range("X4").Style = "Normal"

# Row 5
$ws.Range("A5").Value = 33
$ws.Range("B5").Value = "Feminino"
$ws.Range("C5").Value = "nao"
$ws.Range("D5").Value = "nao"
$ws.Range("E5").Value = "nao"
$ws.Range("F5").Value = "nao"
$ws.Range("G5").Value = "nao"
$ws.Range("H5").Value = "nao"
$ws.Range("I5").Value = "nao"
$ws.Range("J5").Value = "nao"
$ws.Range("K5").Value = "nao"
$ws.Range("L5").Value = "sim"
$ws.Range("M5").Value = "sim"
$ws.Range("N5").Value = "sim"
$ws.Range("O5").Value = "nao"
$ws.Range("P5").Value = "sim"
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 3
$ws.Range("T5").Value = 3
$ws.Range("U5").Value = 3
$ws.Range("V5").Value = "bom"
$ws.Range("W5").Value = "reprovado"
$ws.Range("X5").Value = "'"
$ws.Range("X5").Style = "Normal"

# Row 6
$ws.Range("A6").Value = 33
$ws.Range("B6").Value = "Feminino"
$ws.Range("C6").Value = "nao"
$ws.Range("D6").Value = "nao"
$ws.Range("E6").Value = "nao"
$ws.Range("F6").Value = "nao"
$ws.Range("G6").Value = "nao"
$ws.Range("H6").Value = "nao"
$ws.Range("I6").Value = "nao"
$ws.Range("J6").Value = "nao"
$ws.Range("K6").Value = "nao"
$ws.Range("L6").Value = "sim"
$ws.Range("M6").Value = "sim"
$ws.Range("N6").Value = "sim"
$ws.Range("O6").Value = "nao"
$ws.Range("P6").Value = "sim"
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 3
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = "bom"
$ws.Range("W6").Value = "reprovado"
$ws.Range("X6").Value = "reprovado"

# Row 7
$ws.Range("A7").Value = 33
$ws.Range("B7").Value = "Feminino"
$ws.Range("C7").Value = "nao"
$ws.Range("D7").Value = "nao"
$ws.Range("E7").Value = "nao"
$ws.Range("F7").Value = "nao"
$ws.Range("G7").Value = "nao"
$ws.Range("H7").Value = "nao"
$ws.Range("I7").Value = "nao"
$ws.Range("J7").Value = "nao"
$ws.Range("K7").Value = "nao"
$ws.Range("L7").Value = "sim"
$ws.Range("M7").Value = "sim"
$ws.Range("N7").Value = "sim"
$ws.Range("O7").Value = "nao"
$ws.Range("P7").Value = "sim"
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 1
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 3
$ws.Range("U7").Value = 3
$ws.Range("V7").Value = "bom"
$ws.Range("W7").Value = "reprovado"
$ws.Range("X7").Value = "reprovado"

# Row 8
$ws.Range("A8").Value = 33
$ws.Range("B8").Value = "Feminino"
$ws.Range("C8").Value = "nao"
$ws.Range("D8").Value = "nao"
$ws.Range("E8").Value = "nao"
$ws.Range("F8").Value = "nao"
$ws.Range("G8").Value = "nao"
$ws.Range("H8").Value = "nao"
$ws.Range("I8").Value = "nao"
$ws.Range("J8").Value = "nao"
$ws.Range("K8").Value = "nao"
$ws.Range("L8").Value = "sim"
$ws.Range("M8").Value = "sim"
$ws.Range("N8").Value = "sim"
$ws.Range("O8").Value = "nao"
$ws.Range("P8").Value = "sim"
$ws.Range("Q8").Value = 2
$ws.Range("R8").Value = 1
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 3
$ws.Range("U8").Value = 3
$ws.Range("V8").Value = "bom"
$ws.Range("W8").Value = "reprovado"
$ws.Range("X8").Value = "reprovado"

# Row 9: values kept as text (quote-prefixed numerics) per diff
$ws.Range("A9").Value = "'35"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = "Prefiro não dizer"
$ws.Range("C9").Value = "sim"
$ws.Range("D9").Value = "sim"
$ws.Range("E9").Value = "sim"
$ws.Range("F9").Value = "sim"
$ws.Range("G9").Value = "sim"
$ws.Range("H9").Value = "sim"
$ws.Range("I9").Value = "sim"
$ws.Range("J9").Value = "sim"
$ws.Range("K9").Value = "sim"
$ws.Range("L9").Value = "sim"
$ws.Range("M9").Value = "sim"
$ws.Range("N9").Value = "nao"
$ws.Range("O9").Value = "sim"
$ws.Range("P9").Value = "sim"
$ws.Range("Q9").Value = "'1"
$ws.Range("Q9").Style = "Normal"
$ws.Range("R9").Value = "'1"
$ws.Range("R9").Style = "Normal"
$ws.Range("S9").Value = "'1"
$ws.Range("S9").Style = "Normal"
$ws.Range("T9").Value = "'3"
$ws.Range("T9").Style = "Normal"
$ws.Range("U9").Value = "'2"
$ws.Range("U9").Style = "Normal"
$ws.Range("V9").Value = "ruim"
$ws.Range("W9").Value = "reprovado"
$ws.Range("X9").Value = "reprovado"
